# Append the new "2026/02/06" data row (row 88) to the ModCounts sheet,
# matching the formatting of the existing data rows (e.g. row 87).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 88

# Force column A's date-looking text to stay a literal string instead of
# being auto-parsed into a date serial number (matches how the prior rows
# store their "Date" column as plain text).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2026/02/06"

$ws.Cells.Item($row, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($row, 3).Value = 1175

# Match the existing rows' cell style (center/center alignment).
$ws.Range("A$row`:C$row").Style = $ws.Range("A87:C87").Style
$ws.Range("A$row`:C$row").HorizontalAlignment = -4108
$ws.Range("A$row`:C$row").VerticalAlignment = -4108
